$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value.
# D-column values that look numeric (e.g. "253.99") are prefixed with a
# leading apostrophe so Excel stores them as text, matching the source
# data which keeps all Price/Volume cells as plain strings.
$changes = [ordered]@{
  2  = @{ D = "'42.417.66"; E = '  -0.08%  ' }
  3  = @{ D = "'2.184.85";  E = '  -1.25%  ' }
  4  = @{ E = '  -0.06%  ' }
  5  = @{ D = "'253.99";    E = '  +5.54%  ' }
  6  = @{ D = "'0.613";     E = '  -0.53%  ' }
  7  = @{ D = "'74.03";     E = '  -1.25%  ' }
  8  = @{ E = '  -0.04%  ' }
  9  = @{ E = '  -2.62%  ' }
  10 = @{ D = "'40.57";     E = '  -1.83%  ' }
  11 = @{ E = '  -0.94%  ' }
  12 = @{ B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = "'6.78";  E = '  -1.16%  ' }
  13 = @{ B = 'TRON';     C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = "'0.101"; E = '  -0.36%  ' }
  14 = @{ D = "'2.510.81"; E = '  -1.31%  ' }
  15 = @{ D = "'14.18";    E = '  -3.43%  ' }
  16 = @{ D = "'2.203.45"; E = '  -0.48%  ' }
  17 = @{ D = "'0.769";    E = '  -3.88%  ' }
  18 = @{ D = "'42.326.95";E = '  -0.02%  ' }
  19 = @{ E = '  -2.74%  ' }
  20 = @{ E = '  -0.26%  ' }
  21 = @{ D = "'5.87";     E = '  -0.66%  ' }
  22 = @{ D = "'226.92";   E = '  -0.79%  ' }
  23 = @{ D = "'2.13";     E = '  +1.45%  ' }
  24 = @{ D = "'9.41";     E = '  -6.90%  ' }
  25 = @{ E = '  -0.15%  ' }
  26 = @{ E = '  -4.24%  ' }
  27 = @{ E = '  -0.53%  ' }
  28 = @{ E = '  +1.53%  ' }
  29 = @{ D = "'2.18";     E = '  -2.13%  ' }
  30 = @{ D = "'170.41";   E = '  -1.31%  ' }
  31 = @{ D = "'36.50";    E = '  +8.03%  ' }
  32 = @{ D = "'20.01";    E = '  -0.85%  ' }
  33 = @{ E = '  +1.94%  ' }
  34 = @{ D = "'5.12";     E = '  -4.71%  ' }
  35 = @{ E = '  -0.92%  ' }
  36 = @{ E = '  -0.21%  ' }
  37 = @{ E = '  -3.62%  ' }
  38 = @{ D = "'0.0337";   E = '  +4.53%  ' }
  39 = @{ D = "'11.81";    E = '  -6.45%  ' }
  40 = @{ E = '  -3.15%  ' }
  41 = @{ E = '  +0.23%  ' }
  42 = @{ D = "'59.30";    E = '  -1.98%  ' }
  43 = @{ E = '  -6.19%  ' }
  44 = @{ D = "'102.42";   E = '  +3.02%  ' }
  45 = @{ D = "'0.467";    E = '  +11.12%  ' }
  46 = @{ E = '  +6.25%  ' }
  47 = @{ B = 'Cronos';    C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = "'0.0971"; E = '  -0.82%  ' }
  48 = @{ B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = "'8.25";   E = '  -3.42%  ' }
  49 = @{ E = '  -1.17%  ' }
  50 = @{ E = '  -0.94%  ' }
  51 = @{ E = '  +0.38%  ' }
}

foreach ($row in $changes.Keys) {
  $cols = $changes[$row]
  foreach ($col in $cols.Keys) {
    $ws.Range("$col$row").Value = $cols[$col]
  }
}
